$d = $word.ActiveDocument

# --- Edit 1: "Related work" section / "Irradiance estimation from sky image"
# bullet -> drop the stray duplicated "for" in
# "Camera calibration (adjusting) for based on sun-positions"
$d.Content.Find.Execute(
    "Camera calibration (adjusting) for based on sun-positions", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Camera calibration (adjusting) based on sun-positions", 2)

# --- Edit 2: typo fix "disuse cases" -> "discuss cases" in the bullet about
# showing correlation to diffuse irradiance.
$d.Content.Find.Execute(
    "disuse cases based on images", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "discuss cases based on images", 2)

# --- Word keeps an internal "_GoBack" bookmark marking the last edited spot
# in the document; after the edits above it should now sit right after the
# final text change (inside the "discuss cases" bullet, right before "es").
$r = $d.Content
$found = $r.Find.Execute(
    "discuss cas", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $goBackRange = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
